$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 22.27707562706342
$ws.Cells.Item(2, 3).Value = 3.370399395714863
$ws.Cells.Item(2, 4).Value = 31.33259771631893
$ws.Cells.Item(2, 5).Value = 64.9591273144981
$ws.Cells.Item(2, 6).Value = 27.64850165818148
$ws.Cells.Item(2, 7).Value = 0.7983459015769437
$ws.Cells.Item(2, 8).Value = 23.49999061453126
$ws.Cells.Item(2, 9).Value = 11.76492207051701

$ws.Cells.Item(3, 2).Value = 14.5466383144955
$ws.Cells.Item(3, 3).Value = 7.464427643675682
$ws.Cells.Item(3, 4).Value = 33.72369799056887
$ws.Cells.Item(3, 5).Value = 65.07212332542561
$ws.Cells.Item(3, 6).Value = 25.82769531611721
$ws.Cells.Item(3, 7).Value = 0.7987371938486771
$ws.Cells.Item(3, 8).Value = 22.90503291933921
$ws.Cells.Item(3, 9).Value = 11.91387822353293

$ws.Cells.Item(4, 2).Value = 11.80495454269124
$ws.Cells.Item(4, 3).Value = 10.59418342230815
$ws.Cells.Item(4, 4).Value = 35.37390976589519
$ws.Cells.Item(4, 5).Value = 64.88678185744442
$ws.Cells.Item(4, 6).Value = 24.37960437132714
$ws.Cells.Item(4, 7).Value = 0.7990257511681168
$ws.Cells.Item(4, 8).Value = 22.46651293506611
$ws.Cells.Item(4, 9).Value = 11.7477562833219

$ws.Cells.Item(5, 2).Value = 12.04910518068123
$ws.Cells.Item(5, 3).Value = 11.94843029637465
$ws.Cells.Item(5, 4).Value = 35.91575561381596
$ws.Cells.Item(5, 5).Value = 64.44589251282
$ws.Cells.Item(5, 6).Value = 23.61234108955462
$ws.Cells.Item(5, 7).Value = 0.7992071251969207
$ws.Cells.Item(5, 8).Value = 22.19067861865746
$ws.Cells.Item(5, 9).Value = 11.3139088922775

$ws.Cells.Item(6, 2).Value = 12.36728659535978
$ws.Cells.Item(6, 3).Value = 12.41950439666152
$ws.Cells.Item(6, 4).Value = 36.02443599857974
$ws.Cells.Item(6, 5).Value = 64.12368407801247
$ws.Cells.Item(6, 6).Value = 23.27764444276775
$ws.Cells.Item(6, 7).Value = 0.7993013301232036
$ws.Cells.Item(6, 8).Value = 22.04744434859182
$ws.Cells.Item(6, 9).Value = 10.99373759213834

$ws.Cells.Item(7, 2).Value = 12.36728659535978
$ws.Cells.Item(7, 3).Value = 12.41950439666152
$ws.Cells.Item(7, 4).Value = 36.02443599857974
$ws.Cells.Item(7, 5).Value = 64.12368407801247
$ws.Cells.Item(7, 6).Value = 23.27764444276775
$ws.Cells.Item(7, 7).Value = 0.7993013301232036
$ws.Cells.Item(7, 8).Value = 22.04744434859182
$ws.Cells.Item(7, 9).Value = 10.99373759213834

$ws.Cells.Item(8, 2).Value = 12.18366183725986
$ws.Cells.Item(8, 3).Value = 12.17333238799529
$ws.Cells.Item(8, 4).Value = 35.86034518772338
$ws.Cells.Item(8, 5).Value = 64.06535560894933
$ws.Cells.Item(8, 6).Value = 23.35981922418615
$ws.Cells.Item(8, 7).Value = 0.7992974115562674
$ws.Cells.Item(8, 8).Value = 22.05354620189475
$ws.Cells.Item(8, 9).Value = 10.93407342304822

$ws.Cells.Item(9, 2).Value = 12.15420999185499
$ws.Cells.Item(9, 3).Value = 9.357454706477581
$ws.Cells.Item(9, 4).Value = 34.26001907999796
$ws.Cells.Item(9, 5).Value = 63.98912708626988
$ws.Cells.Item(9, 6).Value = 24.55048402199625
$ws.Cells.Item(9, 7).Value = 0.7991136841700697
$ws.Cells.Item(9, 8).Value = 22.33289532636016
$ws.Cells.Item(9, 9).Value = 10.8423248253931

$ws.Cells.Item(10, 2).Value = 22.60080787853947
$ws.Cells.Item(10, 3).Value = 3.120790783725849
$ws.Cells.Item(10, 4).Value = 30.88485212865321
$ws.Cells.Item(10, 5).Value = 64.33891293853858
$ws.Cells.Item(10, 6).Value = 27.51142198259187
$ws.Cells.Item(10, 7).Value = 0.7984448035169819
$ws.Cells.Item(10, 8).Value = 23.34988826512841
$ws.Cells.Item(10, 9).Value = 11.14317749104012

$ws.Cells.Item(11, 2).Value = 36.61359875409178
$ws.Cells.Item(11, 3).Value = 3.934121698274838
$ws.Cells.Item(11, 4).Value = 27.67034670341199
$ws.Cells.Item(11, 5).Value = 64.61369510743947
$ws.Cells.Item(11, 6).Value = 30.35968418328209
$ws.Cells.Item(11, 7).Value = 0.7976228347915535
$ws.Cells.Item(11, 8).Value = 24.60002757791132
$ws.Cells.Item(11, 9).Value = 11.35321411401351

$ws.Cells.Item(12, 2).Value = 43.36306129878504
$ws.Cells.Item(12, 3).Value = 6.71112887240611
$ws.Cells.Item(12, 4).Value = 25.99016928550224
$ws.Cells.Item(12, 5).Value = 64.55902743687633
$ws.Cells.Item(12, 6).Value = 31.71871937411587
$ws.Cells.Item(12, 7).Value = 0.7971881006981043
$ws.Cells.Item(12, 8).Value = 25.21530439367907
$ws.Cells.Item(12, 9).Value = 11.26100131462551

$ws.Cells.Item(13, 2).Value = 45.29561413767886
$ws.Cells.Item(13, 3).Value = 7.197604302802805
$ws.Cells.Item(13, 4).Value = 24.68161061025798
$ws.Cells.Item(13, 5).Value = 64.48625612947443
$ws.Cells.Item(13, 6).Value = 32.77288289942553
$ws.Cells.Item(13, 7).Value = 0.7968825108552381
$ws.Cells.Item(13, 8).Value = 25.45909213056271
$ws.Cells.Item(13, 9).Value = 11.15818699286311

$ws.Cells.Item(14, 2).Value = 43.11958824529835
$ws.Cells.Item(14, 3).Value = 5.96924025857966
$ws.Cells.Item(14, 4).Value = 24.33270372949186
$ws.Cells.Item(14, 5).Value = 64.34347265647631
$ws.Cells.Item(14, 6).Value = 32.95539153248164
$ws.Cells.Item(14, 7).Value = 0.7968763997880084
$ws.Cells.Item(14, 8).Value = 25.39979554608291
$ws.Cells.Item(14, 9).Value = 11.0102483527408

$ws.Cells.Item(15, 2).Value = 40.76444192238331
$ws.Cells.Item(15, 3).Value = 4.896259637661193
$ws.Cells.Item(15, 4).Value = 24.42410427619701
$ws.Cells.Item(15, 5).Value = 64.11699173644948
$ws.Cells.Item(15, 6).Value = 32.68706780661693
$ws.Cells.Item(15, 7).Value = 0.7970007969636709
$ws.Cells.Item(15, 8).Value = 25.22063427429934
$ws.Cells.Item(15, 9).Value = 10.79053987928659

$ws.Cells.Item(16, 2).Value = 39.89570161153664
$ws.Cells.Item(16, 3).Value = 4.46244646105386
$ws.Cells.Item(16, 4).Value = 24.63039722395499
$ws.Cells.Item(16, 5).Value = 64.32986640505104
$ws.Cells.Item(16, 6).Value = 32.690085996154
$ws.Cells.Item(16, 7).Value = 0.796989571607789
$ws.Cells.Item(16, 8).Value = 25.24514386850154
$ws.Cells.Item(16, 9).Value = 11.00527636592415

$ws.Cells.Item(17, 2).Value = 39.11401342581016
$ws.Cells.Item(17, 3).Value = 4.069904677052261
$ws.Cells.Item(17, 4).Value = 24.88321638243669
$ws.Cells.Item(17, 5).Value = 64.60678793622054
$ws.Cells.Item(17, 6).Value = 32.70712957618439
$ws.Cells.Item(17, 7).Value = 0.7969668350287722
$ws.Cells.Item(17, 8).Value = 25.28650889937001
$ws.Cells.Item(17, 9).Value = 11.28407806775029

$ws.Cells.Item(18, 2).Value = 35.02737491038415
$ws.Cells.Item(18, 3).Value = 2.548354061598416
$ws.Cells.Item(18, 4).Value = 25.84706530169066
$ws.Cells.Item(18, 5).Value = 64.55551118380943
$ws.Cells.Item(18, 6).Value = 31.84341360979079
$ws.Cells.Item(18, 7).Value = 0.7972459558241783
$ws.Cells.Item(18, 8).Value = 24.86387665555799
$ws.Cells.Item(18, 9).Value = 11.25715518253374

$ws.Cells.Item(19, 2).Value = 33.22468995562494
$ws.Cells.Item(19, 3).Value = 2.04333446974085
$ws.Cells.Item(19, 4).Value = 26.58742740623879
$ws.Cells.Item(19, 5).Value = 64.66504675484774
$ws.Cells.Item(19, 6).Value = 31.30988666012703
$ws.Cells.Item(19, 7).Value = 0.7973934058642542
$ws.Cells.Item(19, 8).Value = 24.62695823555432
$ws.Cells.Item(19, 9).Value = 11.38191899645495

$ws.Cells.Item(20, 2).Value = 37.49042325542224
$ws.Cells.Item(20, 3).Value = 3.919602245866761
$ws.Cells.Item(20, 4).Value = 26.5556339609379
$ws.Cells.Item(20, 5).Value = 64.64197317608459
$ws.Cells.Item(20, 6).Value = 31.31521545931107
$ws.Cells.Item(20, 7).Value = 0.7973601396837799
$ws.Cells.Item(20, 8).Value = 24.74384656985828
$ws.Cells.Item(20, 9).Value = 11.35743542684318

$ws.Cells.Item(21, 2).Value = 48.45122925416649
$ws.Cells.Item(21, 3).Value = 8.865953645498765
$ws.Cells.Item(21, 4).Value = 25.07554397554704
$ws.Cells.Item(21, 5).Value = 64.82428710887459
$ws.Cells.Item(21, 6).Value = 32.7156739032883
$ws.Cells.Item(21, 7).Value = 0.7968129441050684
$ws.Cells.Item(21, 8).Value = 25.80375053187469
$ws.Cells.Item(21, 9).Value = 11.49851977318409

$ws.Cells.Item(22, 2).Value = 55.47893207008286
$ws.Cells.Item(22, 3).Value = 11.82400867854242
$ws.Cells.Item(22, 4).Value = 23.84211780243913
$ws.Cells.Item(22, 5).Value = 65.18410086685874
$ws.Cells.Item(22, 6).Value = 34.0812910749634
$ws.Cells.Item(22, 7).Value = 0.7962646198737624
$ws.Cells.Item(22, 8).Value = 26.66374108757664
$ws.Cells.Item(22, 9).Value = 11.81679440829344

$ws.Cells.Item(23, 2).Value = 60.46578833574247
$ws.Cells.Item(23, 3).Value = 13.87173399238003
$ws.Cells.Item(23, 4).Value = 23.01155641305263
$ws.Cells.Item(23, 5).Value = 65.5611069188304
$ws.Cells.Item(23, 6).Value = 35.13153986049915
$ws.Cells.Item(23, 7).Value = 0.7958178756533519
$ws.Cells.Item(23, 8).Value = 27.34345405620727
$ws.Cells.Item(23, 9).Value = 12.16009977088055

$ws.Cells.Item(24, 2).Value = 62.1157881249738
$ws.Cells.Item(24, 3).Value = 14.52608019903177
$ws.Cells.Item(24, 4).Value = 22.81635779503582
$ws.Cells.Item(24, 5).Value = 65.81481545791577
$ws.Cells.Item(24, 6).Value = 35.52333808494438
$ws.Cells.Item(24, 7).Value = 0.7956411603101621
$ws.Cells.Item(24, 8).Value = 27.61254224906727
$ws.Cells.Item(24, 9).Value = 12.40162363983272

$ws.Cells.Item(25, 2).Value = 50.87543064073672
$ws.Cells.Item(25, 3).Value = 9.740112517761961
$ws.Cells.Item(25, 4).Value = 25.25843131601283
$ws.Cells.Item(25, 5).Value = 65.92074962847282
$ws.Cells.Item(25, 6).Value = 33.48320684160136
$ws.Cells.Item(25, 7).Value = 0.7964540688428836
$ws.Cells.Item(25, 8).Value = 26.3764560870963
$ws.Cells.Item(25, 9).Value = 12.57954909938502

